$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, applied with a "@" (text) number format
# so Excel stores them as literal text (matching the source data which used
# inline strings for numeric-looking price / percentage values).
$updates = @{
    'D2' = '309.06'
    'E2' = '-4.25%'
    'D3' = '48.51'
    'E3' = '-3.35%'
    'D4' = '5.185'
    'E4' = '-3.35%'
    'D5' = '0.07752'
    'E5' = '-4.81%'
    'D6' = '4.495'
    'E6' = '-2.30%'
    'D7' = '1.338'
    'E7' = '14.86%'
    'D8' = '1.561'
    'E8' = '-7.06%'
    'D9' = '0.1223'
    'E9' = '-9.24%'
    'D10' = '0.1939'
    'E10' = '-1.68%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.09372'
    'E11' = '-2.32%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.04641'
    'E12' = '1.43%'
    'E13' = '0.02%'
    'D14' = '0.001272'
    'E14' = '-4.29%'
    'D15' = '0.04181'
    'E15' = '-3.08%'
    'D16' = '0.005842'
    'E16' = '-0.34%'
    'D17' = '3.329'
    'E17' = '-1.66%'
    'D18' = '2.274'
    'E18' = '-6.73%'
    'D19' = '0.3491'
    'E19' = '2.85%'
    'D20' = '8.332'
    'E20' = '2.61%'
    'D21' = '0.1339'
    'E21' = '-4.96%'
    'D22' = '0.3039'
    'E22' = '-0.40%'
    'D23' = '0.001276'
    'E23' = '-2.31%'
    'D24' = '0.004178'
    'E24' = '-3.04%'
    'D25' = '0.0001353'
    'E25' = '0.27%'
    'E26' = '-3.99%'
    'D38' = '0.02554'
    'E38' = '-7.59%'
    'D39' = '0.05861'
    'E39' = '5.82%'
    'D40' = '0.01076'
    'E40' = '73.72%'
    'D41' = '0.007949'
    'E41' = '2.54%'
    'D42' = '0.1420'
    'E42' = '-1.82%'
    'D43' = '0.008354'
    'E43' = '8.89%'
    'D44' = '0.007719'
    'E44' = '-12.55%'
    'D45' = '0.3100'
    'E45' = '-11.42%'
    'D46' = '0.00006955'
    'E46' = '2.75%'
    'D47' = '0.00000000752'
    'E47' = '0.27%'
    'D48' = '0.05668'
    'E48' = '-7.57%'
    'E49' = '0.20%'
    'D50' = '0.00002105'
    'E50' = '0.27%'
    'D51' = '0.0002004'
    'E51' = '0.27%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
